$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Black Amber" weekly block (before the
# existing row 29), shifting the old rows 29-40 down to 31-42 and making room
# for a new week's worth of data (2 new rows: Primera + Segunda).
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(29).Insert()

# New row 29: Black Amber / Primera, week of 2022-01-11
$ws.Range("A29").Value = 11
$ws.Range("B29").Value = "Vega Monumental Concepción"
$ws.Range("C29").Value = "Bíobío"
$ws.Range("D29").Value = 44572
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100103
$ws.Range("H29").Value = "Frutos de hueso (carozo)"
$ws.Range("I29").Value = 100103002
$ws.Range("J29").Value = "Ciruela"
$ws.Range("K29").Value = "Black Amber"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 11000
$ws.Range("O29").Value = 12000
$ws.Range("P29").Value = 11500
$ws.Range("Q29").Value = "$/bandeja 18 kilos granel"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 639
$ws.Range("T29").Value = 18

# New row 30: Black Amber / Segunda, same week
$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 44572
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103002
$ws.Range("J30").Value = "Ciruela"
$ws.Range("K30").Value = "Black Amber"
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 10000
$ws.Range("Q30").Value = "$/bandeja 18 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 556
$ws.Range("T30").Value = 18
